$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Asciugamano di carta da bagno kraft riciclato monouso con 1 strato', '0,6409 €', 'Ordine minimo: 30.000 sacchi', 'Dongguan Green Paper Industrial Co., Ltd.', '5.0'),
    @('Fornitura di fabbrica OEM salviette arrotolate piegate usa e getta cotone in microfibra asciugamano per la pulizia del viso', '0,0607 €', 'Ordine minimo: 50.000 sacchi', 'Kunming Ansheng Industry And Trade Co., Ltd.', '4.2'),
    @('OOLIMA usa e getta n piega Z Mulitifold asciugamani di carta per la cucina tovagliolo bagno', '0,0087-0,0174 €', 'Ordine minimo: 100.000 parti', 'Ningbo Oolima Pack Co., Ltd.', '5.0'),
    @('Commercio all''ingrosso di alta qualità 1ply C tovagliolo di carta piegato per 25.4 a mano * 28cm bianco o kraft asciugamano di carta', '0,0867-0,1733 €', 'Ordine minimo: 1.000 sacchi', 'Qingdao Ample Technology Co., Ltd.', $null),
    @('Minguxan Factory Custom Multifold N piega asciugamani di carta bianca marrone carta igienica fazzoletti 1-2ply certificazione FSC Hotel bagno', '7,19-7,97 €', 'Ordine minimo: 500 cartoni', 'Hangzhou Mingxuan Sanitary Products Co., Ltd', '5.0'),
    @('Vendita calda usa e getta Z /C/ V piega a mano asciugamano di carta piegato per fazzoletti di Hotel e ristoranti', '0,2166-0,3118 €', 'Ordine minimo: 10.000 sacchi', 'Zhangzhou Win Share Co., Ltd.', '5.0'),
    @('Molteplici N e V Piega da bagno Asciugamani di Carta A Mano di trasporto del Commercio All''ingrosso Del Tessuto Tovagliolo di Carta', '0,2166-0,563 €', 'Ordine minimo: 20.000 pacchetti', 'Sichuan Petrochemical Yashi Paper Co., Ltd.', '5.0'),
    @('Prezzo competitivo Hotel ristorante fazzoletti usa e getta 1/2 strati riciclati carta vergine N piega mano Inter piegato tovagliolo di carta', '0,4417 €', 'Ordine minimo: 100 scatole', 'Dongguan Xingli Paper Industry Co., Ltd.', $null),
    @('Eco Friendly N fold Z fold Multi fold assorbente usa e getta asciugamano di carta da bagno asciugamano di carta a libro singolo', '0,0693-0,4764 €', 'Ordine minimo: 1 sacco', 'Baoding Clean & Soft Paper Product Technology Co., Ltd.', '3.1'),
    @('Asciugamani per gli ospiti usa e getta all''ingrosso di colore si sentono morbidi asciugamani per le mani da bagno tovaglioli di carta', '0,0607 €', 'Ordine minimo: 5.000 parti', 'Jinan Huaao Nonwoven Technology Co., Ltd.', '4.8'),
    @('Popolare carta bianca monouso da cucina a 2 strati all''ingrosso TravelFold carta da bagno 3 strati', '0,1386-0,1906 €', 'Ordine minimo: 5.000 pacchetti', 'Shandong CT super co.,ltd', '4.9'),
    @('Piccoli asciugamani usa e getta in microfibra arrotolati e piegati per la fabbrica di asciugamani bagnati', '0,026 €', 'Ordine minimo: 1.000 parti', 'Jingzhou Choice Airline Appliance Co., Ltd.', '5.0'),
    @('Tovaglioli 200 bagno monouso asciugamani bianchi per gli ospiti monouso carta velina per le mani a 2 strati piega campione gratuito', '0,0087-0,026 €', 'Ordine minimo: 10.000 parti', 'Nanchang Huayang Trading Co., Ltd.', $null),
    @('Produttore OEM asciugamani di carta ad alta assorbenza personalizzati Asciugamano compattato laminato a 2 strati Multifold V piega Z carta assorbente asciugamano', '0,4417-0,5024 €', 'Ordine minimo: 60.000 sacchi', 'Guizhou Zhiang Trading Co., Ltd.', $null),
    @('Premium usa e getta per la casa asciugamani di carta da cucina spessita asciugamani scatola intera per alberghi e bagni', '0,2252 €', 'Ordine minimo: 5.000 sacchi', 'Baoding Qianzi Trading Co., Ltd.', '4.0'),
    @('Erogatore all''ingrosso di tacchino intrecciato tovagliolo di carta Multifold di alta qualità a basso prezzo 2Ply Z carta velina pieghevole asciugamani per le mani', '5,24-5,37 €', 'Ordine minimo: 2.400 cartoni', 'Erdenay Ticaret Dansmanlk Anonim Sirketi', $null),
    @('Vendita calda pieghevole usa e getta asciugamano e telo da bagno 100% tessuto Non tessuto Extra spesso e assorbente per viaggi Hotel palestra', '0,3811-0,8575 €', 'Ordine minimo: 3.000 parti', 'Jiujiang Yijincheng Medical Equipment Co., Ltd.', '5.0'),
    @('A mano tovagliolo di carta papel higienico V fold articoli per la casa', '0,4946 €', 'Ordine minimo: 20.000 pacchetti', 'Weifang NST Paper Industry Co., Ltd.', '5.0'),
    @('Z salviette di carta a mano pieghevoli e pieghevoli a 1 strato di carta igienica morbida carta per pulizia usa e getta', '0,4764 €', 'Ordine minimo: 30.000 parti', 'Jiangmen Luck Tissue Mfy. Ltd.', $null),
    @('V piega il tovagliolo di carta usa e getta morbida carta igienica pasta di legno asciugamano di carta per le mani tovaglioli dell''hotel', '0,3551 €', 'Ordine minimo: 10 sacchi', 'Baoding Yusen Sanitary Health Supplies Co., Ltd.', '5.0'),
    @('Vendita diretta camera bianca pulizia oem interfold carta tovagliolo asciugamani piegati 3 pieghe carta asciugamani per le mani per il bagno', '8,58-11,18 €', 'Ordine minimo: 200 cartoni', 'Baoding Hozhong Paper Co., Ltd.', $null),
    @('Vendita calda Z Pieghevole di carta alla rinfusa asciugamani bianchi per la pulizia di carta multiforme personalizzati', '8,58-11,26 €', 'Ordine minimo: 200 cartoni', 'Baoding Hozhong Hygienic Products Manufacturing Co., Ltd.', '4.8'),
    @('Qinshe commercio all''ingrosso personalizzato bagno cucina carta igienica tessuto 200 fogli 100% polpa vergine 1 velo N/Z piega asciugamani di carta a mano', '0,5803-0,6323 €', 'Ordine minimo: 2.000 sacchi', 'Hangzhou Qingshe Paper Co., Ltd.', $null),
    @('Salviette umidificate in rotolo di asciugamano singolo asciugamani usa e getta personalizzati all''ingrosso asciugamani per il lavaggio delle mani dell''hotel', '0,1733 €', 'Ordine minimo: 100.000 pacchetti', 'Ningbo Riway Nonwovens Tech Co., Ltd.', '4.7'),
    @('Multi Fold N fold N Fold usa e getta Eco Friendly Premium OEM salviette di carta fabbrica 1 strato di carta', '0,2426 €', 'Ordine minimo: 36.000 sacchi', 'Baoda Paper Enterprise Co., Ltd.', '4.9'),
    @('Asciugamano di carta monouso personalizzabile multiuso 2 3 4 5 C V N Z M', '0,2599 €', 'Ordine minimo: 10.000 pacchetti', 'Weifang Wehicredit Supply Chain Co., Ltd.', '5.0')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = 24 + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt 5; $j++) {
        $cell = $ws.Cells.Item($rowNum, $j + 1)
        $val = $rowVals[$j]
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value = "'" + $val
            $cell.Style = "Normal"
        }
    }
}
